# Append 5 new rows (7-11) to the "Shreyas Gopal " sheet, repeating the
# match-by-match stats for rows 4, 6, 3, 5 and 2 (in that order), exactly
# as they appear in the source diff. All values - including the
# numeric-looking ones in columns G:K - are written as TEXT, matching the
# existing rows 2-6 on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(" Dubai (DSC)", " October 14 2020",    "Capitals won by 13 runs", "Rajasthan Royals", "Delhi Capitals",         "Shreyas Gopal ", "6",  "4",  "1", "0", "150.00"),
    @(" Dubai (DSC)", " September 30 2020",  "KKR won by 37 runs",      "Rajasthan Royals", "Kolkata Knight Riders",  "Shreyas Gopal ", "5",  "7",  "0", "0", "71.42"),
    @(" Dubai (DSC)", " November 01 2020",   "KKR won by 60 runs",      "Rajasthan Royals", "Kolkata Knight Riders",  "Shreyas Gopal ", "23", "23", "2", "0", "100.00"),
    @(" Abu Dhabi",   " October 06 2020",    "Mumbai won by 57 runs",   "Rajasthan Royals", "Mumbai Indians",         "Shreyas Gopal ", "1",  "2",  "0", "0", "50.00"),
    @(" Sharjah",     " October 09 2020",    "Capitals won by 46 runs", "Rajasthan Royals", "Delhi Capitals",         "Shreyas Gopal ", "2",  "3",  "0", "0", "66.66")
)

$startRow = 7
$endRow = $startRow + $newRows.Length - 1
$fillRange = $ws.Range("A$startRow`:K$endRow")

# Force text storage (so "6", "150.00", etc. land as text, same as the
# pre-existing rows) instead of letting Excel auto-convert them to numbers.
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $c = 1 + $j
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}

# Drop the temporary "@" number-format override so the new cells keep the
# workbook's default (unstyled) formatting, matching the other rows.
$fillRange.ClearFormats()

Write-Output "Added rows $startRow to $endRow"
